$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: DataSourceType / numeric / description -------------------------
# A6: copy the "variable name" formatting (blue Consolas, style used by
# RawData/DataSource_noFDR) from A3, then set the text.
$ws.Range("A3").Copy()
$ws.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A6").Value = "DataSourceType"

# B6 already carries the "available values" style (orange Consolas) -
# just fill in the value.
$ws.Range("B6").Value = "numeric"

# C6 uses the default (unstyled) description formatting already present.
$ws.Range("C6").Value = "data source type:  read count / normalized expression values / Fold changes and corrected Pvalue"

# --- Row 7: DataSource_noFDR / True False / description --------------------
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A7").Value = "DataSource_noFDR"

$ws.Range("B7").Value = "True False"

$ws.Range("C7").Value = "Fold-changes only, no corrected P values"
$ws.Range("C7").Font.Name = "Arial"
$ws.Range("C7").Font.Size = 8
$ws.Range("C7").Font.Color = 3355443   # RGB(51,51,51) = #333333

# --- Selection moves to A9 (as last left by the author) --------------------
$null = $ws.Range("A9").Select()
